$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.318.88'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').Value = '1.914.51'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.723'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '255.55'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.72%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.77'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.45%  '

$ws.Range('E9').Value = '  +5.97%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.90'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0766'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0989'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.70%  '

$ws.Range('D13').Value = '2.190.35'
$ws.Range('E13').Value = '  +0.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.727'
$ws.Range('D15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.06%  '

$ws.Range('D17').Value = '1.934.97'
$ws.Range('E17').Value = '  +1.45%  '

$ws.Range('D18').Value = '35.309.25'
$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.66'
$ws.Range('D19').Style = 'Normal'

$ws.Range('D20').Value = '0.0₃0855'
$ws.Range('E20').Value = '  +3.02%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.90%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.38%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.84%  '

$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.52%  '

$ws.Range('E26').Value = '  +4.41%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.71%  '

$ws.Range('E28').Value = '  +2.43%  '

$ws.Range('E29').Value = '  +2.02%  '

$ws.Range('E30').Value = '  +4.45%  '

$ws.Range('D31').Value = '4.129.56'
$ws.Range('E31').Value = '  +19.48%  '

$ws.Range('E32').Value = '  +5.24%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +23.74%  '

$ws.Range('E35').Value = '  +4.10%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.25'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.36%  '

$ws.Range('E37').Value = '  -0.09%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.910'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.02%  '

$ws.Range('E39').Value = '  -0.22%  '

$ws.Range('E40').Value = '  +5.25%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '17.18'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '96.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.80%  '

$ws.Range('E43').Value = '  +1.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0648'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.48%  '

$ws.Range('D45').Value = '1.338.09'
$ws.Range('E45').Value = '  -0.26%  '

$ws.Range('E46').Value = '  +2.05%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.97%  '

$ws.Range('E48').Value = '  +3.12%  '

$ws.Range('E49').Value = '  -0.40%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.64%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.40%  '
